# Adds a "Github Link:" section (with a live hyperlink) after the
# Conclusion's Vernam-Cipher paragraph, inserts a spacer paragraph right
# before the "Conclusion:" heading, and moves the lastRenderedPageBreak
# marker from the "I feel substitution..." run to the "Conclusion:" run
# (Word recomputes this marker whenever pagination shifts because of an
# earlier edit in the same document).

$d = $word.ActiveDocument

function New-FlatOpcPart([string]$bodyInnerXml) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
$bodyInnerXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

# wdParagraph unit used with Range.Expand to snap a found-text range out
# to the whole enclosing paragraph (including its paragraph mark).
$wdParagraph = 4

# ---------------------------------------------------------------------
# 1) Replace the "Conclusion:" paragraph so its first run gains a
#    <w:lastRenderedPageBreak/> before the text. (Its existing
#    <w:bookmarkStart .../> is left alone - the engine keeps it in place
#    across the content replace, so it must not be repeated here.)
# ---------------------------------------------------------------------
$rConclusion = $d.Content
$rConclusion.Find.Execute("Conclusion:") | Out-Null
$rConclusion.Expand($wdParagraph) | Out-Null
$rConclusion.InsertXML((New-FlatOpcPart @'
<w:p w14:paraId="029E33AF" w14:textId="593AD75B" w:rsidR="00AB2044" w:rsidRDefault="002F7378" w:rsidP="009700AA"><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:lastRenderedPageBreak/><w:t>Conclusion:</w:t></w:r><w:r w:rsidR="00AB2044"><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@))

# ---------------------------------------------------------------------
# 2) Replace the "I feel substitution cipher..." paragraph, dropping the
#    <w:lastRenderedPageBreak/> that used to sit in front of its text.
# ---------------------------------------------------------------------
$rIFeel = $d.Content
$rIFeel.Find.Execute("I feel substitution cipher and ROT 13 techniques are") | Out-Null
$rIFeel.Expand($wdParagraph) | Out-Null
$rIFeel.InsertXML((New-FlatOpcPart @'
<w:p w14:paraId="299B45F6" w14:textId="321D8378" w:rsidR="00D521B7" w:rsidRPr="001E4846" w:rsidRDefault="00AB2044" w:rsidP="00267809"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="001E4846"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>I feel substitution cipher and ROT 13 techniques are</w:t></w:r><w:r w:rsidR="00D521B7" w:rsidRPr="001E4846"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="001E4846"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">easier </w:t></w:r><w:r w:rsidR="00D521B7" w:rsidRPr="001E4846"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>to decipher for an attacker using letter frequencies. To make substitution cipher safer and efficient we have to use polyalphabetic cipher technique. Substitution cipher can also be made complex by using permutations and combinations.</w:t></w:r></w:p>
'@))

# ---------------------------------------------------------------------
# 3) Append the new trailing content: a spacer paragraph, a
#    "Github Link:" paragraph, and a paragraph holding the raw URL text
#    (still to be turned into a real hyperlink in step 4).
# ---------------------------------------------------------------------
$endPos = $d.Content.End
$rEnd = $d.Range($endPos, $endPos)
$rEnd.InsertXML((New-FlatOpcPart @'
<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Link:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">https://github.com/Divya-127/CSS-Lab/tree/main/Exp1</w:t></w:r></w:p>
'@))

# ---------------------------------------------------------------------
# 4) Turn the URL text typed above into a live hyperlink (this mints the
#    required external relationship plus the built-in "Hyperlink" /
#    "Unresolved Mention" character styles, exactly as Word does the
#    first time a hyperlink is inserted into a document).
# ---------------------------------------------------------------------
$url = "https://github.com/Divya-127/CSS-Lab/tree/main/Exp1"
$rUrl = $d.Content
$rUrl.Find.Execute($url) | Out-Null
$d.Hyperlinks.Add($rUrl, $url) | Out-Null

# Trailing space run after the hyperlink, matching the source paragraph.
$rAfterLink = $d.Range($rUrl.End, $rUrl.End)
$rAfterLink.InsertAfter(" ")

# ---------------------------------------------------------------------
# 5) Insert the bold spacer paragraph right before "Conclusion:" (done
#    last so it does not disturb the paragraph positions used above).
# ---------------------------------------------------------------------
$rConclusion2 = $d.Content
$rConclusion2.Find.Execute("Conclusion:") | Out-Null
$rConclusion2.Expand($wdParagraph) | Out-Null
$beforePos = $rConclusion2.Start
$rBefore = $d.Range($beforePos, $beforePos)
$rBefore.InsertXML((New-FlatOpcPart '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p>'))

Write-Output "done"
